$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# --- New shared strings must be created in this order so the shared string
# table ends up as: ... 21 = '"case ', 22 = '<-$0.20', 23 = '<-$0.25' ---
$ws.Range("P29").Value = '"case '
$ws.Range("G17").Value = '<-$0.20'
$ws.Range("G21").Value = '<-$0.20'
$ws.Range("G22").Value = '<-$0.25'

# --- Row 17: F17 was the "X" (don't-care) marker; it is corrected to 0,
# with a bold, centered currency style, and an explanatory note in G17 ---
$ws.Range("F17").Value = 0
$ws.Range("F17").NumberFormat = $currencyFmt
$ws.Range("F17").Font.Bold = $true
$ws.Range("F17").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G17").Font.Bold = $true

# --- Row 18: add a bold placeholder cell in G18 (no value) ---
$ws.Range("G18").Font.Bold = $true

# --- Row 19 (new row): bold placeholder cell in G19 (no value) ---
$ws.Range("G19").Font.Bold = $true

# --- Row 20: F20 keeps its value (0.15) but becomes bold currency
# (no special alignment); G20 gets a bold placeholder cell (no value) ---
$ws.Range("F20").NumberFormat = $currencyFmt
$ws.Range("F20").Font.Bold = $true
$ws.Range("G20").Font.Bold = $true

# --- Row 21: F21 corrected from "X" to 0, bold centered currency;
# G21 gets the "<-$0.20" note, bold currency (no special alignment) ---
$ws.Range("F21").Value = 0
$ws.Range("F21").NumberFormat = $currencyFmt
$ws.Range("F21").Font.Bold = $true
$ws.Range("F21").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G21").NumberFormat = $currencyFmt
$ws.Range("G21").Font.Bold = $true

# --- Row 22: F22 corrected from "X" to 0, bold centered currency;
# G22 gets the "<-$0.25" note, bold (no number format) ---
$ws.Range("F22").Value = 0
$ws.Range("F22").NumberFormat = $currencyFmt
$ws.Range("F22").Font.Bold = $true
$ws.Range("F22").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G22").Font.Bold = $true
